$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.963.95'
$ws.Range("E2").Value = '  -3.43%  '

$ws.Range("D3").Value = '2.539.67'
$ws.Range("E3").Value = '  -4.76%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '515.77'
$ws.Range("E5").Value = '  -1.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.76'
$ws.Range("E6").Value = '  -2.49%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -2.16%  '

$ws.Range("D9").Value = '2.544.82'
$ws.Range("E9").Value = '  -5.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.55'
$ws.Range("E10").Value = '  -5.96%  '

$ws.Range("E11").Value = '  -3.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.322'
$ws.Range("E12").Value = '  -3.98%  '

$ws.Range("E13").Value = '  -0.37%  '

$ws.Range("D14").Value = '2.983.82'
$ws.Range("E14").Value = '  -4.86%  '

$ws.Range("D15").Value = '56.959.95'
$ws.Range("E15").Value = '  -3.45%  '

$ws.Range("E16").Value = '  -5.15%  '

$ws.Range("E17").Value = '  -3.44%  '

$ws.Range("D18").Value = '2.521.64'
$ws.Range("E18").Value = '  -6.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '330.00'
$ws.Range("E19").Value = '  -2.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.24'
$ws.Range("E20").Value = '  -3.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.05'
$ws.Range("E21").Value = '  -3.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.13'
$ws.Range("E22").Value = '  -4.04%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.84'
$ws.Range("E24").Value = '  +0.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.167'
$ws.Range("E25").Value = '  +0.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("E27").Value = '  -5.07%  '

$ws.Range("D28").Value = '2.650.70'
$ws.Range("E28").Value = '  -4.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.86'
$ws.Range("E29").Value = '  -4.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").Value = '0.0₃0726'
$ws.Range("E31").Value = '  -9.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.24'
$ws.Range("E32").Value = '  -6.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.56'
$ws.Range("E33").Value = '  -2.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.38'
$ws.Range("E34").Value = '  -1.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.42'
$ws.Range("E35").Value = '  -2.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.96'
$ws.Range("E36").Value = '  -4.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.13'
$ws.Range("E37").Value = '  -5.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.833'
$ws.Range("E38").Value = '  -7.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.56'
$ws.Range("E39").Value = '  -3.65%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.818'
$ws.Range("E40").Value = '  -6.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.42'
$ws.Range("E41").Value = '  -2.63%  '

$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("E43").Value = '  -3.66%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.62'
$ws.Range("E44").Value = '  -0.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '263.47'
$ws.Range("E46").Value = '  -4.46%  '

$ws.Range("E47").Value = '  -6.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.54'
$ws.Range("E48").Value = '  -6.65%  '

$ws.Range("E49").Value = '  -3.65%  '

$ws.Range("D50").Value = '1.951.18'
$ws.Range("E50").Value = '  -4.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0220'
$ws.Range("E51").Value = '  -3.91%  '
